$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.691.30'
$ws.Range('E2').Value = '  +0.21%  '

$ws.Range('D3').Value = '3.506.97'
$ws.Range('E3').Value = '  +0.35%  '

$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.14%  '

$ws.Range('D5').Value = '605.64'
$ws.Range('E5').Value = '  -0.95%  '

$ws.Range('D6').Value = '195.45'
$ws.Range('E6').Value = '  +3.82%  '

$ws.Range('D7').Value = '0.627'
$ws.Range('E7').Value = '  +0.40%  '

$ws.Range('E8').Value = '  +0.14%  '

$ws.Range('D9').Value = '0.202'
$ws.Range('E9').Value = '  -5.55%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.650'
$ws.Range('E10').Value = '  +0.45%  '

$ws.Range('D11').Value = '53.54'
$ws.Range('E11').Value = '  +1.13%  '

$ws.Range('E12').Value = '  -2.08%  '

$ws.Range('D13').Value = '9.51'
$ws.Range('E13').Value = '  +0.11%  '

$ws.Range('D14').Value = '4.063.45'
$ws.Range('E14').Value = '  +0.17%  '

$ws.Range('D15').Value = '595.26'
$ws.Range('E15').Value = '  -0.96%  '

$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = '19.23'
$ws.Range('E16').Value = '  +1.40%  '

$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').Value = '12.81'
$ws.Range('E17').Value = '  +1.80%  '

$ws.Range('D18').Value = '69.823.57'
$ws.Range('E18').Value = '  +0.29%  '

$ws.Range('E19').Value = '  +2.13%  '

$ws.Range('D20').Value = '3.501.95'
$ws.Range('E20').Value = '  +0.96%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.990'
$ws.Range('E21').Value = '  +0.40%  '

$ws.Range('D22').Value = '18.33'
$ws.Range('E22').Value = '  +6.20%  '

$ws.Range('E23').Value = '  +3.54%  '

$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').Value = '4.66'
$ws.Range('E24').Value = '  -0.14%  '

$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '102.02'
$ws.Range('E25').Value = '  -2.90%  '

$ws.Range('D26').Value = '3.17'
$ws.Range('E26').Value = '  +3.98%  '

$ws.Range('E27').Value = '  -0.76%  '

$ws.Range('D28').Value = '9.55'
$ws.Range('E28').Value = '  -1.78%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.30'
$ws.Range('E29').Value = '  -0.13%  '

$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').Value = '7.06'
$ws.Range('E30').Value = '  +1.80%  '

$ws.Range('B31').Value = 'dogwifhat'
$ws.Range('C31').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D31').Value = '4.29'
$ws.Range('E31').Value = '  +3.78%  '

$ws.Range('D32').Value = '12.43'
$ws.Range('E32').Value = '  -0.51%  '

$ws.Range('E33').Value = '  +0.13%  '

$ws.Range('D34').Value = '63.11'
$ws.Range('E34').Value = '  -0.56%  '

$ws.Range('D35').Value = '0.0₃0827'
$ws.Range('E35').Value = '  +6.76%  '

$ws.Range('D36').Value = '3.727.95'
$ws.Range('E36').Value = '  +2.92%  '

$ws.Range('D37').Value = '3.09'
$ws.Range('E37').Value = '  -2.33%  '

$ws.Range('D39').Value = '3.66'
$ws.Range('E39').Value = '  -0.99%  '

$ws.Range('D40').Value = '0.393'
$ws.Range('E40').Value = '  -0.19%  '

$ws.Range('D41').Value = '36.42'
$ws.Range('E41').Value = '  -1.03%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '482.30'
$ws.Range('E42').Value = '  -5.96%  '

$ws.Range('E43').Value = '  -2.74%  '

$ws.Range('D44').Value = '0.0454'
$ws.Range('E44').Value = '  -1.48%  '

$ws.Range('E45').Value = '  -0.97%  '

$ws.Range('D46').Value = '2.82'
$ws.Range('E46').Value = '  -3.24%  '

$ws.Range('E47').Value = '  -1.92%  '

$ws.Range('E48').Value = '  +0.24%  '

$ws.Range('D49').Value = '8.42'
$ws.Range('E49').Value = '  -3.93%  '

$ws.Range('E50').Value = '  +2.37%  '

$ws.Range('E51').Value = '  +10.10%  '
